$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: B1 "C" -> "V", C1 "km" -> "Q1"
$ws.Range("B1").Value = "V"
$ws.Range("C1").Value = "Q1"

# Remove the trailing "V" column (F1); header row shrinks from A1:F1 to A1:E1
$ws.Range("F1").ClearContents()
